# Add a "Save" column (H) to the s_vals sheet, matching the formatting
# already used by the other header cells (e.g. G1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) before writing its text, so it picks up the same bold /
# bordered / centered style used by the rest of row 1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for row 2 - numeric 0, default (unstyled) formatting like
# the other data cells in row 2.
$ws.Range("H2").Value = 0
